$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was recorded for Espinaca. It belongs right
# after the existing row 384, so insert a fresh row at 385 (shifting every
# row below it down by one) and populate it with the new reading.
$ws.Rows.Item(385).Insert()

$ws.Range("A385").Value = 8
$ws.Range("B385").Value = "Terminal La Palmera de La Serena"
$ws.Range("C385").Value = "Coquimbo"
$ws.Range("D385").Value = 45204
$ws.Range("E385").Value = 4
$ws.Range("F385").Value = 100112012
$ws.Range("G385").Value = "Espinaca"
$ws.Range("H385").Value = "Sin especificar"
$ws.Range("I385").Value = "Primera"
$ws.Range("J385").Value = 1800
$ws.Range("K385").Value = 450
$ws.Range("L385").Value = 500
$ws.Range("M385").Value = 475
$ws.Range("N385").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O385").Value = "Provincia del Elquí"
$ws.Range("P385").Value = 950
$ws.Range("Q385").Value = 0.5
$ws.Range("R385").Value = "Hortaliza"
